$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "phone number" style text values for column F (rows 2-18).
# A leading apostrophe forces Excel to store these as text (shared
# strings) rather than numbers, matching the existing convention used
# throughout column F. Resetting the Style back to "Normal" afterwards
# clears the quote-prefix formatting so no new cell style is introduced.
$ws.Range("F2").Value = "'9840013041"
$ws.Range("F3").Value = "'9840001066"
$ws.Range("F4").Value = "'9840020415"
$ws.Range("F5").Value = "'9840061338"
$ws.Range("F6").Value = "'9840027008"
$ws.Range("F7").Value = "'9840073085"
$ws.Range("F8").Value = "'9840055545"
$ws.Range("F9").Value = "'9840034819"
$ws.Range("F10").Value = "'9840078484"
$ws.Range("F11").Value = "'9840019412"
$ws.Range("F12").Value = "'9840082075"
$ws.Range("F13").Value = "'9840054140"
$ws.Range("F14").Value = "'9840054595"
$ws.Range("F15").Value = "'9840011739"
$ws.Range("F16").Value = "'9840025525"
$ws.Range("F17").Value = "'9840095715"
$ws.Range("F18").Value = "'9840015414"
$ws.Range("F2:F18").Style = "Normal"

# AM2 becomes a plain number.
$ws.Range("AM2").Value = 2

# AN2 / AO2 stay text (shared strings).
$ws.Range("AN2").Value = "'3"
$ws.Range("AO2").Value = "'2"
$ws.Range("AN2:AO2").Style = "Normal"

# Update the view/selection to match the committed state.
$ws.Range("AM2").Select()
